$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): drop the "train"/"test" prefixes, keep only one set of metric names
$ws.Range("B1").Value = "mae"
$ws.Range("C1").Value = "mape"
$ws.Range("D1").Value = "rmse"

# Update data rows with the new (revised) values
$ws.Range("B2").Value = 1.23
$ws.Range("C2").Value = 0.04
$ws.Range("D2").Value = 1.7

$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 0.06
$ws.Range("D3").Value = 7.1

# Remove the now-unused "test" columns (E, F, G) entirely so the sheet
# dimension shrinks back down to A1:D3
$ws.Range("E1:G3").Clear()
